$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.542905112095295
$ws.Range("B2").Value = -1.376327953361882

$ws.Range("A3").Value = -0.5187822358978966
$ws.Range("B3").Value = -0.6108791901859103

$ws.Range("A4").Value = -0.9715511871506971
$ws.Range("B4").Value = -0.7840638318940102

$ws.Range("A5").Value = -0.7138522768094833
$ws.Range("B5").Value = -0.6251017597037033

$ws.Range("A6").Value = 0.8231336624746795
$ws.Range("B6").Value = 0.6006606937923294
